# Update the two-digit-divided-by-one-digit practice sheet with a new
# batch of generated division problems. Each cell in the worksheet table
# holds a single "NN÷N=" expression; replace the text of each affected
# cell directly (by row/column) so there is no ambiguity even though
# some old/new values coincide across different cells (e.g. "70÷5=" is
# both a source value in row 1 and a result value in row 17).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of Row -> (Col -> new value), matching the table layout where the
# five problems of each "page row" live on table rows 1, 5, 9, 13, 17.
$updates = @{
    1  = @{ 1 = "92÷7="; 2 = "55÷5="; 3 = "50÷5="; 4 = "37÷7="; 5 = "70÷5=" }
    5  = @{ 1 = "65÷7="; 2 = "74÷6="; 3 = "45÷4="; 4 = "52÷4="; 5 = "26÷3=" }
    9  = @{ 1 = "72÷4="; 2 = "23÷5="; 3 = "63÷4="; 4 = "43÷6="; 5 = "36÷6=" }
    13 = @{ 1 = "13÷8="; 2 = "61÷7="; 3 = "53÷6="; 4 = "62÷6="; 5 = "43÷7=" }
    17 = @{ 1 = "15÷2="; 2 = "51÷6="; 3 = "18÷5="; 4 = "20÷8="; 5 = "75÷2=" }
}

foreach ($rowIndex in $updates.Keys) {
    $cols = $updates[$rowIndex]
    foreach ($colIndex in $cols.Keys) {
        $cell = $t.Cell($rowIndex, $colIndex)
        $cell.Range.Text = $cols[$colIndex]
    }
}
